$wb = $excel.ActiveWorkbook

# Delete the "Desarquivamentos Pendentes" worksheet entirely.
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()

# Rename "Paineis DARQ" -> "PAINEIS DARQ" (uppercase).
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"

# Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO" (uppercase, accented).
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Keep the first tab active/selected (deleting a later sheet can shift focus).
$wb.Worksheets.Item("PAINEIS DARQ").Activate()
